# Generate Report for handoff
# Marks the handoff as failed for both localized sheets: clears the stale
# "Latest Handoff File" hyperlink/value, resets the handoff datetime, flips
# the status message, and marks the handoff reason as Ignored.

$wb = $excel.ActiveWorkbook

# The Overview sheet's per-locale status column mirrors the same status text
# as the detail sheets below, so it also flips to the failure message.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status (B2): "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Remove the hyperlink that lived on C2 ("Latest Handoff File") and
    # clear the cell entirely so it drops out of the sheet.
    $linksToRemove = @()
    foreach ($link in $ws.Hyperlinks) {
        if ($link.Range.Address() -eq '$C$2') {
            $linksToRemove += $link
        }
    }
    foreach ($link in $linksToRemove) {
        $link.Delete()
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime (D2): reset to the zero-value timestamp.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (H2): "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
